$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers (losing formatting like trailing zeros)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row-by-row cell value updates matching the refreshed crypto data
$ws.Range("D2").Value = '30.638.82'
$ws.Range("E2").Value = '  +2.45%  '
$ws.Range("D3").Value = '1.676.24'
$ws.Range("E3").Value = '  +2.71%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '219.22'
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("D6").Value = '0.528'
$ws.Range("E6").Value = '  +1.87%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '29.02'
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("E9").Value = '  +1.92%  '
$ws.Range("E10").Value = '  +5.63%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '1.918.05'
$ws.Range("E12").Value = '  +2.78%  '
$ws.Range("D13").Value = '1.687.36'
$ws.Range("E13").Value = '  +3.30%  '
$ws.Range("D14").Value = '10.09'
$ws.Range("E14").Value = '  +8.32%  '
$ws.Range("D15").Value = '0.604'
$ws.Range("E15").Value = '  +7.33%  '
$ws.Range("E16").Value = '  +4.22%  '
$ws.Range("D17").Value = '30.649.07'
$ws.Range("D18").Value = '65.90'
$ws.Range("E18").Value = '  +2.78%  '
$ws.Range("D19").Value = '242.69'
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("D20").Value = '0.0₃0720'
$ws.Range("E20").Value = '  +2.56%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").Value = '4.23'
$ws.Range("E22").Value = '  +2.41%  '
$ws.Range("D23").Value = '9.95'
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("D25").Value = '159.29'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").Value = '15.79'
$ws.Range("E26").Value = '  +2.16%  '
$ws.Range("E27").Value = '  +2.35%  '
$ws.Range("D28").Value = '6.68'
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").Value = '1.15'
$ws.Range("E31").Value = '  +3.99%  '
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("E33").Value = '  +4.08%  '
$ws.Range("D34").Value = '1.509.91'
$ws.Range("E34").Value = '  +6.08%  '
$ws.Range("D35").Value = '1.77'
$ws.Range("E35").Value = '  +6.16%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '1.02'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("B37").Value = 'Aave'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D37").Value = '83.09'
$ws.Range("E37").Value = '  +10.06%  '
$ws.Range("D38").Value = '0.600'
$ws.Range("E38").Value = '  +7.96%  '
$ws.Range("E39").Value = '  +4.38%  '
$ws.Range("E40").Value = '  -2.91%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '2.01'
$ws.Range("E42").Value = '  +1.45%  '
$ws.Range("D43").Value = '0.836'
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("D44").Value = '0.0497'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  +4.19%  '
$ws.Range("D48").Value = '1.810.87'
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("D49").Value = '49.72'
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0116'
$ws.Range("E50").Value = '  +3.41%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '92.77'
$ws.Range("E51").Value = '  +2.54%  '
